$d = $word.ActiveDocument

# --- Fix 1: "dia" -> "día" (typo correction) ---------------------------
# The word "dia" sits in its own run, flanked by w:proofErr spellStart/spellEnd
# markers (Word's automatic spell-check annotations). Replacing across the
# whole "hoy dia." span merges the touched runs and drops the now-irrelevant
# proofErr markers, exactly as happens when Word re-paginates/re-checks the
# text after a correction.
$d.Content.Find.Execute("hoy dia.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "hoy día.", 2)

# Re-split the merged run back into the original three runs (the text before
# "día", "día" itself, and the trailing period) by toggling a character
# property on just the "día" sub-range - this forces Word to give it its own
# run again without reintroducing the proofErr markers.
$rng = $d.Range(0, $d.Content.End)
$rng.Find.Execute("día", $false, $true, $false, $false, $false, $true, 1, `
    $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

# --- Fix 2: merge the "con" split caused by a grammar-check annotation -
# "todo con lo que cuenta" was split into three runs around the word "con"
# because of w:proofErr gramStart/gramEnd markers. Re-asserting the same
# text via Find/Replace merges the runs back into one contiguous run and
# drops the stale proofErr markers.
$d.Content.Find.Execute("todo con lo que cuenta", $false, $false, $false, `
    $false, $false, $true, 1, $false, "todo con lo que cuenta", 2)
